$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# SARAALERT-1260: allow the vaccine table to be populated on import.
#
# Adds ten new "vaccine" columns (CY:DH) to the Monitorees export sheet:
#   Vaccine 1 Group Name / Product Name / Administration Date / Dose Number / Notes
#   Vaccine 2 Group Name / Product Name / Administration Date / Dose Number / Notes
# and fills in sample vaccination data for the first few monitorees.
# ---------------------------------------------------------------------------

# ---- Header row (row 1) ---------------------------------------------------
$ws.Range("CY1").Value = "Vaccine 1 Group Name"
$ws.Range("CZ1").Value = "Vaccine 1 Product Name"
$ws.Range("DA1").Value = "Vaccine 1 Administration Date"
$ws.Range("DB1").Value = "Vaccine 1 Dose Number"
$ws.Range("DD1").Value = "Vaccine 2 Group Name"
$ws.Range("DE1").Value = "Vaccine 2 Product Name"
$ws.Range("DF1").Value = "Vaccine 2 Administration Date"
$ws.Range("DG1").Value = "Vaccine 2 Dose Number"
$ws.Range("DH1").Value = "Vaccine 2 Notes"
$ws.Range("DC1").Value = "Vaccine 1 Notes"

# Administration-date columns should be stored as text (not auto-converted
# to a date serial number), matching the rest of the free-text date columns
# already present in the workbook (e.g. columns BR/BS). Only format the
# cells that will actually receive a value.
$ws.Range("DA1:DA7").NumberFormat = "@"
$ws.Range("DF1:DF4").NumberFormat = "@"

# ---- Monitoree 1 (row 2) - fully populated, two vaccine doses -------------
$ws.Range("CY2").Value = "COVID-19"
$ws.Range("CZ2").Value = "Moderna COVID-19 Vaccine"
$ws.Range("DA2").Value = "2020-06-01"
$ws.Range("DF2").Value = "2020-06-20"
$ws.Range("DC2").Value = "notes 1"
$ws.Range("DH2").Value = "notes 2"
$ws.Range("DB2").Value = 1
$ws.Range("DG2").Value = 2

# ---- Vaccine product/group names for the remaining monitorees -------------
$ws.Range("CZ6").Value = "Janssen (J&J) COVID-19 Vaccine"
$ws.Range("CZ3").Value = "Pfizer-BioNTech COVID-19 Vaccine"
$ws.Range("CZ4").Value = "Unknown"
$ws.Range("CZ7").Value = "Unknown"
$ws.Range("CZ5").Value = "Moderna COVID-19 Vaccine"

$ws.Range("CY3").Value = "COVID-19"
$ws.Range("CY4").Value = "COVID-19"
$ws.Range("CY5").Value = "COVID-19"
$ws.Range("CY6").Value = "COVID-19"
$ws.Range("CY7").Value = "COVID-19"

# ---- Administration dates (vaccine 1) for the remaining monitorees --------
$ws.Range("DA3").Value = "2020-06-02"
$ws.Range("DA6").Value = "2020-06-03"
$ws.Range("DA4").Value = "2020-06-04"
$ws.Range("DA7").Value = "2020-06-02"
$ws.Range("DA5").Value = "2020-06-01"

$ws.Range("DB3").Value = 1
$ws.Range("DB4").Value = 1
$ws.Range("DB5").Value = 1
$ws.Range("DB6").Value = 1
$ws.Range("DB7").Value = 1

# ---- Administration dates (vaccine 2) for monitorees with a second dose ---
$ws.Range("DF3").Value = "2020-06-21"
$ws.Range("DF4").Value = "2020-06-22"

$ws.Range("DG3").Value = 2
$ws.Range("DG4").Value = 2

# ---- Normalized vaccine enum columns (group/product name, vaccine 2) ------
$ws.Range("DD2").Value = "covid19"
$ws.Range("DE2").Value = "moderna covid19 vaccine"

$ws.Range("DD3").Value = "COVID-19"
$ws.Range("DE3").Value = "pfizerbiontech covid19 vaccine"
$ws.Range("DD4").Value = "COVID-19"
$ws.Range("DE4").Value = "unknown"

# ---------------------------------------------------------------------------
# Column widths - approximate the autosized widths Excel would compute for
# the newly-populated columns.
# ---------------------------------------------------------------------------
$ws.Range("CY1:DH7").EntireColumn.AutoFit()

# Reset the view back to the top-left of the sheet / A1 selected.
$ws.Range("A1").Select()
